$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto prices/volumes (refreshed data pull)
$ws.Range("D2").Value = "44.222.76"
$ws.Range("E2").Value = "  +6.13%  "

$ws.Range("D3").Value = "2.283.26"
$ws.Range("E3").Value = "  +3.70%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.639"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.96"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +7.95%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.434"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +8.27%  "

$ws.Range("E10").Value = "  +17.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.78"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.01"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +17.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.19%  "

$ws.Range("D14").Value = "2.624.89"
$ws.Range("E14").Value = "  +3.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.95"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.96%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.03"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +8.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.834"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.15%  "

$ws.Range("D18").Value = "2.283.87"
$ws.Range("E18").Value = "  +3.69%  "

$ws.Range("D19").Value = "44.089.09"
$ws.Range("E19").Value = "  +5.94%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000103"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +14.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.16"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.26"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.55%  "

$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("E25").Value = "  +6.05%  "

$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.12"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.78"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.49"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.92%  "

$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("D31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.83"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.61%  "

$ws.Range("E33").Value = "  +3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0688"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.83%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.75"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.99"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.11%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.88"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +9.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.84"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +8.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.36"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.04%  "

$ws.Range("E40").Value = "  +4.07%  "

$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.55"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.75"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +9.67%  "

$ws.Range("E44").Value = "  +2.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.13%  "

$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.47"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("D48").Value = "1.473.40"
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("E49").Value = "  +6.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000208"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -12.50%  "

$ws.Range("E51").Value = "  +2.16%  "

